# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx data refresh described in the commit
# "Updated cryptos list on Fri Jun 23 15:51:20 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.817.03"
$ws.Range("E2").Value = "  +3.49%  "
$ws.Range("D3").Value = "1.897.42"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.71"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9984"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4934"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2955"
$ws.Range("E8").Value = "  +2.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06693"
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("D10").Value = "1.889.43"
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.76"
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07243"
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6767"
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.049"
$ws.Range("E14").Value = "  +5.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.97"
$ws.Range("E15").Value = "  +3.19%  "
$ws.Range("D16").Value = "30.711.05"
$ws.Range("E16").Value = "  +3.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007909"
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9976"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.00"
$ws.Range("E19").Value = "  +2.10%  "
$ws.Range("D20").Value = "2.128.41"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.818"
$ws.Range("E22").Value = "  +1.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.956"
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "161.55"
$ws.Range("E24").Value = "  +20.62%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.268"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.80"
$ws.Range("E26").Value = "  +3.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.36"
$ws.Range("E27").Value = "  +4.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.927"
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.427"
$ws.Range("E29").Value = "  +3.66%  "
$ws.Range("E30").Value = "  +2.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08855"
$ws.Range("E31").Value = "  +2.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.012"
$ws.Range("E32").Value = "  +2.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05241"
$ws.Range("E33").Value = "  +3.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7391"
$ws.Range("E34").Value = "  +4.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.128"
$ws.Range("E35").Value = "  +2.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.662"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01849"
$ws.Range("E37").Value = "  +13.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.709"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.195"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9414"
$ws.Range("E40").Value = "  +1.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.855"
$ws.Range("E41").Value = "  -3.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4290"
$ws.Range("E42").Value = "  +3.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.53"
$ws.Range("E43").Value = "  +2.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.573"
$ws.Range("E45").Value = "  +0.87%  "
$ws.Range("E46").Value = "  +4.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05797"
$ws.Range("E47").Value = "  +1.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "33.00"
$ws.Range("E48").Value = "  +1.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.401"
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3817"
$ws.Range("E50").Value = "  +3.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.359"
$ws.Range("E51").Value = "  +1.59%  "
